# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update Home ("H") row target-depth tallies ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 448
$wsOff.Range("C2").Value = 309
$wsOff.Range("D2").Value = 112
$wsOff.Range("E2").Value = 67
$wsOff.Range("F2").Value = 7

# --- DEF sheet: update Home ("H") row target-depth tallies ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 300
$wsDef.Range("C2").Value = 200
$wsDef.Range("D2").Value = 78
$wsDef.Range("E2").Value = 43
